$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.682.93"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.530.74"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "621.49"
$ws.Range("E5").Value = "  +4.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.37"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "3.527.98"
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("E11").Value = "  -2.80%  "
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.20"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000275"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "4.100.21"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.45"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "607.56"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").Value = "3.536.92"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "70.841.78"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("E21").Value = "  +2.12%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.08"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.59"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.55"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -1.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.64"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.10"
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.11"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.29"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("E34").Value = "  -4.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "616.42"
$ws.Range("E35").Value = "  -8.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0494"
$ws.Range("E36").Value = "  +3.82%  "
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.38"
$ws.Range("E41").Value = "  -6.41%  "
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").Value = "3.349.79"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").Value = "0.0₃0726"
$ws.Range("E44").Value = "  +3.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.310"
$ws.Range("E45").Value = "  -2.70%  "
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "31.84"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.49"
$ws.Range("E48").Value = "  -4.47%  "
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.00"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("E51").Value = "  +7.42%  "
